# Apply "Add data for 2021-10-11" update to the carjacking-by-neighborhood-by-month workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet and update the running "through" date label (row 1, column B)
$ws.Name = "Through 2021-10-03"
$ws.Cells.Item(1, 2).Value = "October 2021 (through October 03)"

# Cell value updates for the new day's counts (row, column letter, new value)
$ws.Range("L3").Value = 1

$ws.Range("B4").Value = 3
$ws.Range("L4").Value = 3

$ws.Range("B5").Value = 1

$ws.Range("B6").Value = 4

$ws.Range("B8").Value = 2

$ws.Range("L10").Value = 1

$ws.Range("B12").Value = 1
$ws.Range("AF12").Value = 1

$ws.Range("AF13").Value = 1

$ws.Range("B15").Value = 3

$ws.Range("B20").Value = 2

$ws.Range("L37").Value = 2
$ws.Range("V37").Value = 1
$ws.Range("AP37").Value = 2

$ws.Range("B65").Value = 1

$ws.Range("B78").Value = 1

$ws.Range("B87").Value = 1
$ws.Range("AF87").Value = 1

$ws.Range("B89").Value = 1
